# "Generate Report for Handoff"
#
# The b.md file needed a new handoff round (a newer source revision exists),
# so its status flips from "Handed back: in sync with en-US" to
# "Ready for handoff", and the per-language sheets get fresh handoff file
# names / timestamps plus an explanatory error detail message.

$wb = $excel.ActiveWorkbook

$status_ready = "Ready for handoff"
$overviewDate = "2016-09-04 10:41:47"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8b1657ef765d6d84f7d453234bd3d51fba8cd890/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1c32cb42516adcbc5d58f6cf72fc3984e64528e0/e2e/b.md."

# ---------------------------------------------------------------
# Overview sheet: row 3 is the b.md entry.
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $status_ready
$wsOverview.Range("F3").Value = $status_ready
$wsOverview.Range("G3").Value = $overviewDate

# ---------------------------------------------------------------
# zh-cn sheet: row 3 is the b.md entry.
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $status_ready
# Leading "'" forces text (instead of boolean) storage for "False"; the
# subsequent Style reset drops the quote-prefix cell style so the cell
# ends up a plain shared-string reference, same as its neighbours.
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-04 10:41:42"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1640625

# ---------------------------------------------------------------
# de-de sheet: row 3 is the b.md entry.
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $status_ready
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = $overviewDate
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1640625
